$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title paragraph: collapse the three "Blue Waters "/"Petascale"/
#    " Semester Curriculum v1.0" runs (split apart for spell-check markers)
#    into one continuous run with the same text (and drop the proofErr tags).
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Blue Waters Petascale Semester Curriculum v1.0", $true, $false, $false, $false, $false, `
              $true, 1, $false, "Blue Waters Petascale Semester Curriculum v1.0", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Swap the licence from CC BY-NC 4.0 to CC BY-SA 4.0 (text + hyperlink).
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("CC BY-NC", $true, $false, $false, $false, $false, `
              $true, 1, $false, "CC BY-SA", 2) | Out-Null

foreach ($h in $d.Hyperlinks) {
    if ($h.Address -eq "https://creativecommons.org/licenses/by-nc/4.0") {
        $h.TextToDisplay = "https://creativecommons.org/licenses/by-sa/4.0"
        $h.Address = "https://creativecommons.org/licenses/by-sa/4.0"
    }
}

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of the document (left over from
#    the last edit position) up to the very start of the document, alongside
#    the other bookmarks already anchored there.
# ---------------------------------------------------------------------------
foreach ($b in $d.Bookmarks) {
    if ($b.Name -eq "_GoBack") {
        $b.Delete()
    }
}
$startRange = $d.Range(0, 0)
$d.Bookmarks.Add("_GoBack", $startRange) | Out-Null

# ---------------------------------------------------------------------------
# 4. Drop the stale cached page-break marker in front of "Read the " — Word
#    recalculates this on layout, so re-stamping the run's text clears it.
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Read the", $true, $false, $false, $false, $false, `
              $true, 1, $false, "Read the", 2) | Out-Null
